# Updated symbol list: refresh Price (D) and Volume(1h) (E) figures for the
# affected rows. Values are stored as literal text in this sheet (not
# numbers/percentages), so each write is entered with a leading apostrophe —
# exactly like a user typing '261.05 into the cell — to force text entry and
# keep the General number format instead of letting Excel auto-convert the
# numeric-looking / percent-looking text into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.05"
$ws.Range("E2").Value = "'1.29%"
$ws.Range("D3").Value = "'27.15"
$ws.Range("E3").Value = "'1.75%"
$ws.Range("D4").Value = "'4.670"
$ws.Range("E4").Value = "'0.67%"
$ws.Range("D5").Value = "'0.06179"
$ws.Range("E5").Value = "'3.91%"
$ws.Range("E6").Value = "'0.67%"
$ws.Range("D7").Value = "'0.8513"
$ws.Range("E7").Value = "'-0.67%"
$ws.Range("D8").Value = "'0.9180"
$ws.Range("E8").Value = "'0.45%"
$ws.Range("D9").Value = "'0.1413"
$ws.Range("E9").Value = "'1.85%"
$ws.Range("D10").Value = "'0.04886"
$ws.Range("E10").Value = "'10.46%"
$ws.Range("D11").Value = "'0.07087"
$ws.Range("E11").Value = "'0.95%"
$ws.Range("D12").Value = "'0.03107"
$ws.Range("E12").Value = "'2.99%"
$ws.Range("D13").Value = "'0.09039"
$ws.Range("E13").Value = "'-0.92%"
$ws.Range("D14").Value = "'0.001549"
$ws.Range("E14").Value = "'1.57%"
$ws.Range("D15").Value = "'0.0006171"
$ws.Range("E15").Value = "'2.57%"
$ws.Range("D16").Value = "'0.006085"
$ws.Range("E16").Value = "'0.05%"
$ws.Range("E17").Value = "'-0.65%"
$ws.Range("E18").Value = "'0.79%"
$ws.Range("E19").Value = "'0.23%"
$ws.Range("E20").Value = "'-0.81%"
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'0.28%"
$ws.Range("D22").Value = "'4.102"
$ws.Range("E22").Value = "'6.47%"
$ws.Range("D23").Value = "'0.04245"
$ws.Range("E23").Value = "'0.92%"
$ws.Range("E24").Value = "'0.18%"
$ws.Range("E25").Value = "'-15.50%"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("E27").Value = "'-8.00%"
$ws.Range("E40").Value = "'1.48%"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.58%"
$ws.Range("D42").Value = "'0.004103"
$ws.Range("E42").Value = "'9.11%"
$ws.Range("E43").Value = "'9.88%"
$ws.Range("E44").Value = "'-5.01%"
$ws.Range("D45").Value = "'0.00005165"
$ws.Range("E45").Value = "'1.40%"
$ws.Range("E46").Value = "'0.20%"
$ws.Range("E47").Value = "'8.22%"
$ws.Range("D48").Value = "'0.1624"
$ws.Range("E48").Value = "'-32.59%"
$ws.Range("E49").Value = "'0.20%"
$ws.Range("E50").Value = "'0.20%"
